$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Remove the old hidden "_GoBack" bookmark while it is still unambiguous
# (there is exactly one in the document before we start editing). Word
# always keeps this bookmark pinned to wherever the user last edited, so
# once we finish typing the new paragraphs below we will re-create it in
# its new location.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# Insert the three new paragraphs right after the "Set up" paragraph that
# ends with "...loading screens. " (paragraph 6), describing the key
# bindings for the game.
# ---------------------------------------------------------------------
$pSetup = $d.Paragraphs.Item(6)

$pSetup.Range.InsertParagraphAfter()
$pMove = $d.Paragraphs.Item(7)
$pMove.Range.Text = "WASD and arrow keys to move"

$pMove.Range.InsertParagraphAfter()
$pSpell = $d.Paragraphs.Item(8)

# This paragraph contains the "_GoBack" bookmark sitting right after the
# "Q", splitting the paragraph's text into two runs. Build it via raw
# WordOpenXML so the run boundaries / xml:space flags come out exactly as
# intended instead of being merged & re-normalised by the Text setter.
$spellXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Q</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/><w:r><w:t xml:space="preserve">, E are spell changing keys, well left and right mouse are firing. </w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$pSpell.Range.InsertXML($spellXml)

$pSpell.Range.InsertParagraphAfter()
$pInteract = $d.Paragraphs.Item(9)
$pInteract.Range.Text = "F is the interact button for things like levers, and other buttons. "
